# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E20) listed the years in descending
# order (2006..2002). This updates it to ascending order (2002..2006),
# which is how the refreshed account-statement database now reports it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2002"
$ws.Range("E17").Value = "2003"
$ws.Range("E18").Value = "2004"
$ws.Range("E19").Value = "2005"
$ws.Range("E20").Value = "2006"
